$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3 ; $ws.Range("G2").Value = 39.860211 ; $ws.Range("H2").Value = 119.580633 ; $ws.Range("I2").Value = 0.08465532840534135 ; $ws.Range("J2").Value = 0.08465532840534136 ; $ws.Range("K2").Value = 3 ; $ws.Range("M2").Value = 4.471321666666667 ; $ws.Range("N2").Value = 13.413965 ; $ws.Range("O2").Value = 0.083204941376588 ; $ws.Range("P2").Value = 0.08320494137658797 ; $ws.Range("Q2").Value = 178.227825082205 ; $ws.Range("R2").Value = 1604.050425739845 ; $ws.Range("S2").Value = 0.007043741637182232 ; $ws.Range("T2").Value = 0.007043741637182231
$ws.Range("E3").Value = 3 ; $ws.Range("G3").Value = 39.860211 ; $ws.Range("H3").Value = 119.580633 ; $ws.Range("I3").Value = 0.08465532840534135 ; $ws.Range("J3").Value = 0.08465532840534136 ; $ws.Range("K3").Value = 3 ; $ws.Range("M3").Value = 0.4515893333333333 ; $ws.Range("N3").Value = 1.354768 ; $ws.Range("O3").Value = 0.008403435674603098 ; $ws.Range("P3").Value = 0.008403435674603096 ; $ws.Range("Q3").Value = 18.000446112016 ; $ws.Range("R3").Value = 162.004015008144 ; $ws.Range("S3").Value = 0.0007113956067666865 ; $ws.Range("T3").Value = 0.0007113956067666865
$ws.Range("E4").Value = 3 ; $ws.Range("G4").Value = 39.860211 ; $ws.Range("H4").Value = 119.580633 ; $ws.Range("I4").Value = 0.08465532840534135 ; $ws.Range("J4").Value = 0.08465532840534136 ; $ws.Range("K4").Value = 3 ; $ws.Range("M4").Value = 6.212987666666666 ; $ws.Range("N4").Value = 18.638963 ; $ws.Range("O4").Value = 0.1156148703038507 ; $ws.Range("P4").Value = 0.1156148703038506 ; $ws.Range("Q4").Value = 247.650999333731 ; $ws.Range("R4").Value = 2228.858994003579 ; $ws.Range("S4").Value = 0.009787414814113425 ; $ws.Range("T4").Value = 0.009787414814113425
$ws.Range("E5").Value = 3 ; $ws.Range("G5").Value = 39.860211 ; $ws.Range("H5").Value = 119.580633 ; $ws.Range("I5").Value = 0.08465532840534135 ; $ws.Range("J5").Value = 0.08465532840534136 ; $ws.Range("K5").Value = 3 ; $ws.Range("M5").Value = 6.473970333333334 ; $ws.Range("N5").Value = 19.421911 ; $ws.Range("O5").Value = 0.1204713868104106 ; $ws.Range("P5").Value = 0.1204713868104106 ; $ws.Range("Q5").Value = 258.053823494407 ; $ws.Range("R5").Value = 2322.484411449663 ; $ws.Range("S5").Value = 0.01019854481388222 ; $ws.Range("T5").Value = 0.01019854481388222
$ws.Range("E6").Value = 3 ; $ws.Range("G6").Value = 39.860211 ; $ws.Range("H6").Value = 119.580633 ; $ws.Range("I6").Value = 0.08465532840534135 ; $ws.Range("J6").Value = 0.08465532840534136 ; $ws.Range("K6").Value = 3 ; $ws.Range("M6").Value = 32.24961033333333 ; $ws.Range("N6").Value = 96.748831 ; $ws.Range("O6").Value = 0.6001194137310196 ; $ws.Range("P6").Value = 0.6001194137310196 ; $ws.Range("Q6").Value = 1285.476272554447 ; $ws.Range("R6").Value = 11569.28645299002 ; $ws.Range("S6").Value = 0.05080330605182038 ; $ws.Range("T6").Value = 0.05080330605182039
$ws.Range("E7").Value = 3 ; $ws.Range("G7").Value = 39.860211 ; $ws.Range("H7").Value = 119.580633 ; $ws.Range("I7").Value = 0.08465532840534135 ; $ws.Range("J7").Value = 0.08465532840534136 ; $ws.Range("K7").Value = 3 ; $ws.Range("M7").Value = 3.879176 ; $ws.Range("N7").Value = 11.637528 ; $ws.Range("O7").Value = 0.07218595210352802 ; $ws.Range("P7").Value = 0.072185952103528 ; $ws.Range("Q7").Value = 154.624773866136 ; $ws.Range("R7").Value = 1391.622964795224 ; $ws.Range("S7").Value = 0.006110925481576405 ; $ws.Range("T7").Value = 0.006110925481576405
$ws.Range("E8").Value = 3 ; $ws.Range("G8").Value = 41.8492 ; $ws.Range("H8").Value = 125.5476 ; $ws.Range("I8").Value = 0.08887955383630085 ; $ws.Range("J8").Value = 0.08887955383630086 ; $ws.Range("K8").Value = 3 ; $ws.Range("M8").Value = 4.471321666666667 ; $ws.Range("N8").Value = 13.413965 ; $ws.Range("O8").Value = 0.083204941376588 ; $ws.Range("P8").Value = 0.08320494137658797 ; $ws.Range("Q8").Value = 187.1212346926667 ; $ws.Range("R8").Value = 1684.091112234 ; $ws.Range("S8").Value = 0.007395218066526709 ; $ws.Range("T8").Value = 0.007395218066526708
$ws.Range("E9").Value = 3 ; $ws.Range("G9").Value = 41.8492 ; $ws.Range("H9").Value = 125.5476 ; $ws.Range("I9").Value = 0.08887955383630085 ; $ws.Range("J9").Value = 0.08887955383630086 ; $ws.Range("K9").Value = 3 ; $ws.Range("M9").Value = 0.4515893333333333 ; $ws.Range("N9").Value = 1.354768 ; $ws.Range("O9").Value = 0.008403435674603098 ; $ws.Range("P9").Value = 0.008403435674603096 ; $ws.Range("Q9").Value = 18.89865232853333 ; $ws.Range("R9").Value = 170.0878709568 ; $ws.Range("S9").Value = 0.0007468936134507771 ; $ws.Range("T9").Value = 0.0007468936134507771
$ws.Range("E10").Value = 3 ; $ws.Range("G10").Value = 41.8492 ; $ws.Range("H10").Value = 125.5476 ; $ws.Range("I10").Value = 0.08887955383630085 ; $ws.Range("J10").Value = 0.08887955383630086 ; $ws.Range("K10").Value = 3 ; $ws.Range("M10").Value = 6.212987666666666 ; $ws.Range("N10").Value = 18.638963 ; $ws.Range("O10").Value = 0.1156148703038507 ; $ws.Range("P10").Value = 0.1156148703038506 ; $ws.Range("Q10").Value = 260.0085634598667 ; $ws.Range("R10").Value = 2340.0770711388 ; $ws.Range("S10").Value = 0.01027579808944804 ; $ws.Range("T10").Value = 0.01027579808944803
$ws.Range("E11").Value = 3 ; $ws.Range("G11").Value = 41.8492 ; $ws.Range("H11").Value = 125.5476 ; $ws.Range("I11").Value = 0.08887955383630085 ; $ws.Range("J11").Value = 0.08887955383630086 ; $ws.Range("K11").Value = 3 ; $ws.Range("M11").Value = 6.473970333333334 ; $ws.Range("N11").Value = 19.421911 ; $ws.Range("O11").Value = 0.1204713868104106 ; $ws.Range("P11").Value = 0.1204713868104106 ; $ws.Range("Q11").Value = 270.9304792737334 ; $ws.Range("R11").Value = 2438.374313463601 ; $ws.Range("S11").Value = 0.01070744310974971 ; $ws.Range("T11").Value = 0.01070744310974971
$ws.Range("E12").Value = 3 ; $ws.Range("G12").Value = 41.8492 ; $ws.Range("H12").Value = 125.5476 ; $ws.Range("I12").Value = 0.08887955383630085 ; $ws.Range("J12").Value = 0.08887955383630086 ; $ws.Range("K12").Value = 3 ; $ws.Range("M12").Value = 32.24961033333333 ; $ws.Range("N12").Value = 96.748831 ; $ws.Range("O12").Value = 0.6001194137310196 ; $ws.Range("P12").Value = 0.6001194137310196 ; $ws.Range("Q12").Value = 1349.620392761733 ; $ws.Range("R12").Value = 12146.5835348556 ; $ws.Range("S12").Value = 0.05333834574091546 ; $ws.Range("T12").Value = 0.05333834574091547
$ws.Range("E13").Value = 3 ; $ws.Range("G13").Value = 41.8492 ; $ws.Range("H13").Value = 125.5476 ; $ws.Range("I13").Value = 0.08887955383630085 ; $ws.Range("J13").Value = 0.08887955383630086 ; $ws.Range("K13").Value = 3 ; $ws.Range("M13").Value = 3.879176 ; $ws.Range("N13").Value = 11.637528 ; $ws.Range("O13").Value = 0.07218595210352802 ; $ws.Range("P13").Value = 0.072185952103528 ; $ws.Range("Q13").Value = 162.3404122592 ; $ws.Range("R13").Value = 1461.0637103328 ; $ws.Range("S13").Value = 0.006415855216210153 ; $ws.Range("T13").Value = 0.006415855216210153
$ws.Range("E14").Value = 3 ; $ws.Range("G14").Value = 47.89520899999999 ; $ws.Range("H14").Value = 143.685627 ; $ws.Range("I14").Value = 0.1017200999497333 ; $ws.Range("J14").Value = 0.1017200999497333 ; $ws.Range("K14").Value = 3 ; $ws.Range("M14").Value = 4.471321666666667 ; $ws.Range("N14").Value = 13.413965 ; $ws.Range("O14").Value = 0.083204941376588 ; $ws.Range("P14").Value = 0.08320494137658797 ; $ws.Range("Q14").Value = 214.1548857312283 ; $ws.Range("R14").Value = 1927.393971581055 ; $ws.Range("S14").Value = 0.008463614953138235 ; $ws.Range("T14").Value = 0.008463614953138233
$ws.Range("E15").Value = 3 ; $ws.Range("G15").Value = 47.89520899999999 ; $ws.Range("H15").Value = 143.685627 ; $ws.Range("I15").Value = 0.1017200999497333 ; $ws.Range("J15").Value = 0.1017200999497333 ; $ws.Range("K15").Value = 3 ; $ws.Range("M15").Value = 0.4515893333333333 ; $ws.Range("N15").Value = 1.354768 ; $ws.Range("O15").Value = 0.008403435674603098 ; $ws.Range("P15").Value = 0.008403435674603096 ; $ws.Range("Q15").Value = 21.62896550217067 ; $ws.Range("R15").Value = 194.660689519536 ; $ws.Range("S15").Value = 0.0008547983167417818 ; $ws.Range("T15").Value = 0.0008547983167417818
$ws.Range("E16").Value = 3 ; $ws.Range("G16").Value = 47.89520899999999 ; $ws.Range("H16").Value = 143.685627 ; $ws.Range("I16").Value = 0.1017200999497333 ; $ws.Range("J16").Value = 0.1017200999497333 ; $ws.Range("K16").Value = 3 ; $ws.Range("M16").Value = 6.212987666666666 ; $ws.Range("N16").Value = 18.638963 ; $ws.Range("O16").Value = 0.1156148703038507 ; $ws.Range("P16").Value = 0.1156148703038506 ; $ws.Range("Q16").Value = 297.5723428094223 ; $ws.Range("R16").Value = 2678.1510852848 ; $ws.Range("S16").Value = 0.01176035616298315 ; $ws.Range("T16").Value = 0.01176035616298315
$ws.Range("E17").Value = 3 ; $ws.Range("G17").Value = 47.89520899999999 ; $ws.Range("H17").Value = 143.685627 ; $ws.Range("I17").Value = 0.1017200999497333 ; $ws.Range("J17").Value = 0.1017200999497333 ; $ws.Range("K17").Value = 3 ; $ws.Range("M17").Value = 6.473970333333334 ; $ws.Range("N17").Value = 19.421911 ; $ws.Range("O17").Value = 0.1204713868104106 ; $ws.Range("P17").Value = 0.1204713868104106 ; $ws.Range("Q17").Value = 310.0721621747996 ; $ws.Range("R17").Value = 2790.649459573197 ; $ws.Range("S17").Value = 0.01225436150743795 ; $ws.Range("T17").Value = 0.01225436150743795
$ws.Range("E18").Value = 3 ; $ws.Range("G18").Value = 47.89520899999999 ; $ws.Range("H18").Value = 143.685627 ; $ws.Range("I18").Value = 0.1017200999497333 ; $ws.Range("J18").Value = 0.1017200999497333 ; $ws.Range("K18").Value = 3 ; $ws.Range("M18").Value = 32.24961033333333 ; $ws.Range("N18").Value = 96.748831 ; $ws.Range("O18").Value = 0.6001194137310196 ; $ws.Range("P18").Value = 0.6001194137310196 ; $ws.Range("Q18").Value = 1544.601827083559 ; $ws.Range("R18").Value = 13901.41644375203 ; $ws.Range("S18").Value = 0.06104420674649468 ; $ws.Range("T18").Value = 0.06104420674649469
$ws.Range("E19").Value = 3 ; $ws.Range("G19").Value = 47.89520899999999 ; $ws.Range("H19").Value = 143.685627 ; $ws.Range("I19").Value = 0.1017200999497333 ; $ws.Range("J19").Value = 0.1017200999497333 ; $ws.Range("K19").Value = 3 ; $ws.Range("M19").Value = 3.879176 ; $ws.Range("N19").Value = 11.637528 ; $ws.Range("O19").Value = 0.07218595210352802 ; $ws.Range("P19").Value = 0.072185952103528 ; $ws.Range("Q19").Value = 185.793945267784 ; $ws.Range("R19").Value = 1672.145507410056 ; $ws.Range("S19").Value = 0.007342762262937533 ; $ws.Range("T19").Value = 0.007342762262937532
$ws.Range("E20").Value = 3 ; $ws.Range("G20").Value = 59.96530133333332 ; $ws.Range("H20").Value = 179.895904 ; $ws.Range("I20").Value = 0.1273546263289621 ; $ws.Range("J20").Value = 0.1273546263289621 ; $ws.Range("K20").Value = 3 ; $ws.Range("M20").Value = 4.471321666666667 ; $ws.Range("N20").Value = 13.413965 ; $ws.Range("O20").Value = 0.083204941376588 ; $ws.Range("P20").Value = 0.08320494137658797 ; $ws.Range("Q20").Value = 268.1241510999289 ; $ws.Range("R20").Value = 2413.11735989936 ; $ws.Range("S20").Value = 0.01059653421773856 ; $ws.Range("T20").Value = 0.01059653421773856
$ws.Range("E21").Value = 3 ; $ws.Range("G21").Value = 59.96530133333332 ; $ws.Range("H21").Value = 179.895904 ; $ws.Range("I21").Value = 0.1273546263289621 ; $ws.Range("J21").Value = 0.1273546263289621 ; $ws.Range("K21").Value = 3 ; $ws.Range("M21").Value = 0.4515893333333333 ; $ws.Range("N21").Value = 1.354768 ; $ws.Range("O21").Value = 0.008403435674603098 ; $ws.Range("P21").Value = 0.008403435674603096 ; $ws.Range("Q21").Value = 27.07969045225244 ; $ws.Range("R21").Value = 243.717214070272 ; $ws.Range("S21").Value = 0.001070216410218547 ; $ws.Range("T21").Value = 0.001070216410218547
$ws.Range("E22").Value = 3 ; $ws.Range("G22").Value = 59.96530133333332 ; $ws.Range("H22").Value = 179.895904 ; $ws.Range("I22").Value = 0.1273546263289621 ; $ws.Range("J22").Value = 0.1273546263289621 ; $ws.Range("K22").Value = 3 ; $ws.Range("M22").Value = 6.212987666666666 ; $ws.Range("N22").Value = 18.638963 ; $ws.Range("O22").Value = 0.1156148703038507 ; $ws.Range("P22").Value = 0.1156148703038506 ; $ws.Range("Q22").Value = 372.5636776119501 ; $ws.Range("R22").Value = 3353.073098507551 ; $ws.Range("S22").Value = 0.01472408860561832 ; $ws.Range("T22").Value = 0.01472408860561832
$ws.Range("E23").Value = 3 ; $ws.Range("G23").Value = 59.96530133333332 ; $ws.Range("H23").Value = 179.895904 ; $ws.Range("I23").Value = 0.1273546263289621 ; $ws.Range("J23").Value = 0.1273546263289621 ; $ws.Range("K23").Value = 3 ; $ws.Range("M23").Value = 6.473970333333334 ; $ws.Range("N23").Value = 19.421911 ; $ws.Range("O23").Value = 0.1204713868104106 ; $ws.Range("P23").Value = 0.1204713868104106 ; $ws.Range("Q23").Value = 388.2135818613937 ; $ws.Range("R23").Value = 3493.922236752544 ; $ws.Range("S23").Value = 0.0153425884505717 ; $ws.Range("T23").Value = 0.0153425884505717
$ws.Range("E24").Value = 3 ; $ws.Range("G24").Value = 59.96530133333332 ; $ws.Range("H24").Value = 179.895904 ; $ws.Range("I24").Value = 0.1273546263289621 ; $ws.Range("J24").Value = 0.1273546263289621 ; $ws.Range("K24").Value = 3 ; $ws.Range("M24").Value = 32.24961033333333 ; $ws.Range("N24").Value = 96.748831 ; $ws.Range("O24").Value = 0.6001194137310196 ; $ws.Range("P24").Value = 0.6001194137310196 ; $ws.Range("Q24").Value = 1933.857601520913 ; $ws.Range("R24").Value = 17404.71841368822 ; $ws.Range("S24").Value = 0.07642798368846983 ; $ws.Range("T24").Value = 0.07642798368846983
$ws.Range("E25").Value = 3 ; $ws.Range("G25").Value = 59.96530133333332 ; $ws.Range("H25").Value = 179.895904 ; $ws.Range("I25").Value = 0.1273546263289621 ; $ws.Range("J25").Value = 0.1273546263289621 ; $ws.Range("K25").Value = 3 ; $ws.Range("M25").Value = 3.879176 ; $ws.Range("N25").Value = 11.637528 ; $ws.Range("O25").Value = 0.07218595210352802 ; $ws.Range("P25").Value = 0.072185952103528 ; $ws.Range("Q25").Value = 232.6159577650346 ; $ws.Range("R25").Value = 2093.543619885312 ; $ws.Range("S25").Value = 0.009193214956345168 ; $ws.Range("T25").Value = 0.009193214956345166
$ws.Range("E26").Value = 3 ; $ws.Range("G26").Value = 243.1290336666667 ; $ws.Range("H26").Value = 729.387101 ; $ws.Range("I26").Value = 0.5163587365336566 ; $ws.Range("J26").Value = 0.5163587365336566 ; $ws.Range("K26").Value = 3 ; $ws.Range("M26").Value = 4.471321666666667 ; $ws.Range("N26").Value = 13.413965 ; $ws.Range("O26").Value = 0.083204941376588 ; $ws.Range("P26").Value = 0.08320494137658797 ; $ws.Range("Q26").Value = 1087.108116029496 ; $ws.Range("R26").Value = 9783.973044265465 ; $ws.Range("S26").Value = 0.04296359840257195 ; $ws.Range("T26").Value = 0.04296359840257193
$ws.Range("E27").Value = 3 ; $ws.Range("G27").Value = 243.1290336666667 ; $ws.Range("H27").Value = 729.387101 ; $ws.Range("I27").Value = 0.5163587365336566 ; $ws.Range("J27").Value = 0.5163587365336566 ; $ws.Range("K27").Value = 3 ; $ws.Range("M27").Value = 0.4515893333333333 ; $ws.Range("N27").Value = 1.354768 ; $ws.Range("O27").Value = 0.008403435674603098 ; $ws.Range("P27").Value = 0.008403435674603096 ; $ws.Range("Q27").Value = 109.7944782275076 ; $ws.Range("R27").Value = 988.1503040475681 ; $ws.Range("S27").Value = 0.004339187427479912 ; $ws.Range("T27").Value = 0.004339187427479911
$ws.Range("E28").Value = 3 ; $ws.Range("G28").Value = 243.1290336666667 ; $ws.Range("H28").Value = 729.387101 ; $ws.Range("I28").Value = 0.5163587365336566 ; $ws.Range("J28").Value = 0.5163587365336566 ; $ws.Range("K28").Value = 3 ; $ws.Range("M28").Value = 6.212987666666666 ; $ws.Range("N28").Value = 18.638963 ; $ws.Range("O28").Value = 0.1156148703038507 ; $ws.Range("P28").Value = 0.1156148703038506 ; $ws.Range("Q28").Value = 1510.557687579585 ; $ws.Range("R28").Value = 13595.01918821626 ; $ws.Range("S28").Value = 0.05969874835459891 ; $ws.Range("T28").Value = 0.05969874835459889
$ws.Range("E29").Value = 3 ; $ws.Range("G29").Value = 243.1290336666667 ; $ws.Range("H29").Value = 729.387101 ; $ws.Range("I29").Value = 0.5163587365336566 ; $ws.Range("J29").Value = 0.5163587365336566 ; $ws.Range("K29").Value = 3 ; $ws.Range("M29").Value = 6.473970333333334 ; $ws.Range("N29").Value = 19.421911 ; $ws.Range("O29").Value = 0.1204713868104106 ; $ws.Range("P29").Value = 0.1204713868104106 ; $ws.Range("Q29").Value = 1574.010151130002 ; $ws.Range("R29").Value = 14166.09136017001 ; $ws.Range("S29").Value = 0.06220645308188104 ; $ws.Range("T29").Value = 0.06220645308188103
$ws.Range("E30").Value = 3 ; $ws.Range("G30").Value = 243.1290336666667 ; $ws.Range("H30").Value = 729.387101 ; $ws.Range("I30").Value = 0.5163587365336566 ; $ws.Range("J30").Value = 0.5163587365336566 ; $ws.Range("K30").Value = 3 ; $ws.Range("M30").Value = 32.24961033333333 ; $ws.Range("N30").Value = 96.748831 ; $ws.Range("O30").Value = 0.6001194137310196 ; $ws.Range("P30").Value = 0.6001194137310196 ; $ws.Range("Q30").Value = 7840.816596469881 ; $ws.Range("R30").Value = 70567.34936822893 ; $ws.Range("S30").Value = 0.309876902243468 ; $ws.Range("T30").Value = 0.309876902243468
$ws.Range("E31").Value = 3 ; $ws.Range("G31").Value = 243.1290336666667 ; $ws.Range("H31").Value = 729.387101 ; $ws.Range("I31").Value = 0.5163587365336566 ; $ws.Range("J31").Value = 0.5163587365336566 ; $ws.Range("K31").Value = 3 ; $ws.Range("M31").Value = 3.879176 ; $ws.Range("N31").Value = 11.637528 ; $ws.Range("O31").Value = 0.07218595210352802 ; $ws.Range("P31").Value = 0.072185952103528 ; $ws.Range("Q31").Value = 943.1403123029254 ; $ws.Range("R31").Value = 8488.262810726328 ; $ws.Range("S31").Value = 0.03727384702365678 ; $ws.Range("T31").Value = 0.03727384702365677
$ws.Range("E32").Value = 3 ; $ws.Range("G32").Value = 38.15399366666667 ; $ws.Range("H32").Value = 114.461981 ; $ws.Range("I32").Value = 0.0810316549460057 ; $ws.Range("J32").Value = 0.08103165494600571 ; $ws.Range("K32").Value = 3 ; $ws.Range("M32").Value = 4.471321666666667 ; $ws.Range("N32").Value = 13.413965 ; $ws.Range("O32").Value = 0.083204941376588 ; $ws.Range("P32").Value = 0.08320494137658797 ; $ws.Range("Q32").Value = 170.5987785516295 ; $ws.Range("R32").Value = 1535.389006964665 ; $ws.Range("S32").Value = 0.006742234099430311 ; $ws.Range("T32").Value = 0.00674223409943031
$ws.Range("E33").Value = 3 ; $ws.Range("G33").Value = 38.15399366666667 ; $ws.Range("H33").Value = 114.461981 ; $ws.Range("I33").Value = 0.0810316549460057 ; $ws.Range("J33").Value = 0.08103165494600571 ; $ws.Range("K33").Value = 3 ; $ws.Range("M33").Value = 0.4515893333333333 ; $ws.Range("N33").Value = 1.354768 ; $ws.Range("O33").Value = 0.008403435674603098 ; $ws.Range("P33").Value = 0.008403435674603096 ; $ws.Range("Q33").Value = 17.22993656393422 ; $ws.Range("R33").Value = 155.069429075408 ; $ws.Range("S33").Value = 0.0006809442999453929 ; $ws.Range("T33").Value = 0.0006809442999453928
$ws.Range("E34").Value = 3 ; $ws.Range("G34").Value = 38.15399366666667 ; $ws.Range("H34").Value = 114.461981 ; $ws.Range("I34").Value = 0.0810316549460057 ; $ws.Range("J34").Value = 0.08103165494600571 ; $ws.Range("K34").Value = 3 ; $ws.Range("M34").Value = 6.212987666666666 ; $ws.Range("N34").Value = 18.638963 ; $ws.Range("O34").Value = 0.1156148703038507 ; $ws.Range("P34").Value = 0.1156148703038506 ; $ws.Range("Q34").Value = 237.0502920850781 ; $ws.Range("R34").Value = 2133.452628765703 ; $ws.Range("S34").Value = 0.009368464277088829 ; $ws.Range("T34").Value = 0.009368464277088828
$ws.Range("E35").Value = 3 ; $ws.Range("G35").Value = 38.15399366666667 ; $ws.Range("H35").Value = 114.461981 ; $ws.Range("I35").Value = 0.0810316549460057 ; $ws.Range("J35").Value = 0.08103165494600571 ; $ws.Range("K35").Value = 3 ; $ws.Range("M35").Value = 6.473970333333334 ; $ws.Range("N35").Value = 19.421911 ; $ws.Range("O35").Value = 0.1204713868104106 ; $ws.Range("P35").Value = 0.1204713868104106 ; $ws.Range("Q35").Value = 247.0078230961879 ; $ws.Range("R35").Value = 2223.070407865691 ; $ws.Range("S35").Value = 0.009761995846887974 ; $ws.Range("T35").Value = 0.009761995846887974
$ws.Range("E36").Value = 3 ; $ws.Range("G36").Value = 38.15399366666667 ; $ws.Range("H36").Value = 114.461981 ; $ws.Range("I36").Value = 0.0810316549460057 ; $ws.Range("J36").Value = 0.08103165494600571 ; $ws.Range("K36").Value = 3 ; $ws.Range("M36").Value = 32.24961033333333 ; $ws.Range("N36").Value = 96.748831 ; $ws.Range("O36").Value = 0.6001194137310196 ; $ws.Range("P36").Value = 0.6001194137310196 ; $ws.Range("Q36").Value = 1230.451428410468 ; $ws.Range("R36").Value = 11074.06285569421 ; $ws.Range("S36").Value = 0.04862866925985122 ; $ws.Range("T36").Value = 0.04862866925985122
$ws.Range("E37").Value = 3 ; $ws.Range("G37").Value = 38.15399366666667 ; $ws.Range("H37").Value = 114.461981 ; $ws.Range("I37").Value = 0.0810316549460057 ; $ws.Range("J37").Value = 0.08103165494600571 ; $ws.Range("K37").Value = 3 ; $ws.Range("M37").Value = 3.879176 ; $ws.Range("N37").Value = 11.637528 ; $ws.Range("O37").Value = 0.07218595210352802 ; $ws.Range("P37").Value = 0.072185952103528 ; $ws.Range("Q37").Value = 148.0060565358853 ; $ws.Range("R37").Value = 1332.054508822968 ; $ws.Range("S37").Value = 0.005849347162801976 ; $ws.Range("T37").Value = 0.005849347162801976
